# Generate Report for Handoff
#
# The localization-status report tracks each source file's handoff state
# per target locale. Once the handoff (HO) xliff files for
# "ecb07ba0-6028-4097-95e1-7a0eb13f3d0b.md" were generated, its status
# moves from "In Translation" to "Ready for handoff" for both the zh-cn
# and de-de locales, its priority becomes machine-translation ("mt"
# instead of "ht"), and the "Latest Handoff Datetime" / rollup
# "Latest HO Xliff Generate Date" timestamps are refreshed to reflect the
# new handoff generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for ecb07ba0-6028-4097-95e1-7a0eb13f3d0b.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-22 12:12:37"

# --- zh-cn sheet: row for ecb07ba0-6028-4097-95e1-7a0eb13f3d0b.md ------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-22 12:12:33"

# --- de-de sheet: row for ecb07ba0-6028-4097-95e1-7a0eb13f3d0b.md ------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-22 12:12:37"

# --- Cosmetic: the longer "Ready for handoff" status text widened the
# Status columns in Excel's autofit. Nudge the affected columns to the
# closest width this engine can represent.
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
